# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The "Periodo Mora" (column E) listing for the worker's overdue periods is
# reordered from descending (2501, 2412, 2411, 2410, 2409, 2408, 2407) to
# ascending (2407, 2408, 2409, 2410, 2411, 2412, 2501), and the "Valor Mora"
# (column F) amount that belongs to period 2501 (50266) travels with it from
# row 16 down to row 22, while the remaining rows keep the 52000 figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2407"
$ws.Range("E17").Value = "2408"
$ws.Range("E18").Value = "2409"
$ws.Range("E19").Value = "2410"
$ws.Range("E20").Value = "2411"
$ws.Range("E21").Value = "2412"
$ws.Range("E22").Value = "2501"

$ws.Range("F16").Value = 52000
$ws.Range("F17").Value = 52000
$ws.Range("F18").Value = 52000
$ws.Range("F19").Value = 52000
$ws.Range("F20").Value = 52000
$ws.Range("F21").Value = 52000
$ws.Range("F22").Value = 50266
